$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43, shifting existing rows 43:118 down to 44:119
$ws.Rows("43").Insert()

# Populate the newly inserted row 43 with the new record
$ws.Range("A43").Value2 = 11
$ws.Range("B43").Value2 = "Vega Monumental Concepción"
$ws.Range("C43").Value2 = "Bíobío"
$ws.Range("D43").Value2 = 45044
$ws.Range("E43").Value2 = 8
$ws.Range("F43").Value2 = 100112012
$ws.Range("G43").Value2 = "Espinaca"
$ws.Range("H43").Value2 = "Sin especificar"
$ws.Range("I43").Value2 = "Primera"
$ws.Range("J43").Value2 = 220
$ws.Range("K43").Value2 = 7000
$ws.Range("L43").Value2 = 8000
$ws.Range("M43").Value2 = 7545
$ws.Range("N43").Value2 = "$/cuna 10 kilos"
$ws.Range("O43").Value2 = "Región Metropolitana"
$ws.Range("P43").Value2 = 754
$ws.Range("Q43").Value2 = 10
$ws.Range("R43").Value2 = "Hortaliza"
